$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (old Fonttype column shifts right to K)
$ws.Columns("J:J").Insert()

# Match the new column's width to the neighboring "File" column (I)
$ws.Columns("J:J").ColumnWidth = $ws.Columns("I:I").ColumnWidth

# New "Centering" column header + sample data
$ws.Range("J1").Value2 = "Centering"
$ws.Range("J2").Value2 = "x"
$ws.Range("J3").Value2 = "yx"

# Update the example-image text in B3 to the new multi-line copy, with wrap text enabled
$ws.Range("B3").WrapText = $true
$ws.Range("B3").Value2 = "An Example Image `nsdfgfhgjfhg,j.h,ghghfg`ntest"

# Update the active selection shown when the workbook is opened
$ws.Range("B3").Select()
